# Adjust schedule data for the "difference instance": Flight No. (B) becomes
# a plain numeric id matching column A, and Arrival/Departure times (C/D),
# Location (E) and Type (F) are updated to the new instance's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    # row, B,  C,                     D,                     E,     F
    @(2,  1,  0.02361111111111111,  0.05069444444444444,  "E4",  3),
    @(3,  2,  0.03611111111111111,  0.06180555555555556,  "C26", 1),
    @(4,  3,  0.004166666666666667, 0.02569444444444444,  "D49", 1),
    @(5,  4,  0.006944444444444444, 0.02777777777777778,  "A14", 3),
    @(6,  5,  0.01736111111111111,  0.04444444444444445,  "D40", 1),
    @(7,  6,  0.01180555555555556,  0.03263888888888889,  "E28", 3),
    @(8,  7,  0.005555555555555556, 0.02847222222222222,  "D37", 3),
    @(9,  8,  0.009722222222222222, 0.03402777777777777,  "A15", 3),
    @(10, 9,  0.03680555555555556,  0.05902777777777778,  "E8",  3),
    @(11, 10, 0.03263888888888889,  0.05555555555555555,  "B4",  3),
    @(12, 11, 0.03958333333333333,  0.06666666666666667,  "F42", 1),
    @(13, 12, 0.03888888888888889,  0.06180555555555556,  "A3",  2),
    @(14, 13, 0.04027777777777778,  0.06458333333333334,  "B10", 1),
    @(15, 14, 0.03055555555555555,  0.05555555555555555,  "C24", 1),
    @(16, 15, 0.01041666666666667,  0.03611111111111111,  "E7",  1),
    @(17, 16, 0.01805555555555555,  0.04027777777777778,  "F59", 1),
    @(18, 17, 0.03541666666666667,  0.06180555555555556,  "F30", 1),
    @(19, 18, 0.05694444444444444,  0.0798611111111111,   "B3",  1),
    @(20, 19, 0.07777777777777778,  0.1034722222222222,   "D48", 3),
    @(21, 20, 0.04791666666666667,  0.07013888888888889,  "D46", 3)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
}
